# unitTest_base_macro2.xlsx — add `storeKeys(json,jsonpath,var)` expression
# and drop the stray `text` entry from the `target` category list (and its
# orphaned single-cell `text` data column) on the hidden "#system" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("#system")

# ------------------------------------------------------------------
# 1) json: insert "storeKeys(json,jsonpath,var)" as the new M16, pushing
#    the existing storeValue/storeValues rows down one cell within
#    column M only (column 13). Using plain cell-by-cell value copy
#    (rather than Range.Insert, which shifts the *whole row* in this
#    host) keeps columns A, E, G, K, etc. on row 16-18 untouched.
# ------------------------------------------------------------------
$oldM17 = $ws.Cells.Item(17, 13).Text
$oldM16 = $ws.Cells.Item(16, 13).Text
$ws.Cells.Item(18, 13).Value = $oldM17
$ws.Cells.Item(17, 13).Value = $oldM16
$ws.Cells.Item(16, 13).Value = "storeKeys(json,jsonpath,var)"

# ------------------------------------------------------------------
# 2) target: remove the "text" entry (row 25 in column A / column 1),
#    shifting webalert/webcookie/ws/ws.async/xml/etc. up one cell
#    within column A only, and clearing the now-stale last row (A31).
# ------------------------------------------------------------------
for ($r = 25; $r -le 30; $r++) {
    $below = $ws.Cells.Item($r + 1, 1).Text
    $ws.Cells.Item($r, 1).Value = $below
}
$ws.Cells.Item(31, 1).ClearContents()

# ------------------------------------------------------------------
# 3) Remove the orphaned single-cell "text" column entirely (column Y),
#    shifting web/webalert/webcookie/ws/ws.async/xml left one column.
# ------------------------------------------------------------------
$ws.Columns.Item(25).Delete()

# ------------------------------------------------------------------
# 4) Update the defined names (workbook-level) to match the new
#    locations/extents of the data they point to.
# ------------------------------------------------------------------
$wb.Names.Item("json").RefersTo = "='#system'!`$M`$2:`$M`$18"
$wb.Names.Item("target").RefersTo = "='#system'!`$A`$2:`$A`$30"
$wb.Names.Item("web").RefersTo = "='#system'!`$Y`$2:`$Y`$129"
$wb.Names.Item("webalert").RefersTo = "='#system'!`$Z`$2:`$Z`$8"
$wb.Names.Item("webcookie").RefersTo = "='#system'!`$AA`$2:`$AA`$8"
$wb.Names.Item("ws").RefersTo = "='#system'!`$AB`$2:`$AB`$17"
$wb.Names.Item("ws.async").RefersTo = "='#system'!`$AC`$2:`$AC`$8"
$wb.Names.Item("xml").RefersTo = "='#system'!`$AD`$2:`$AD`$27"

Write-Host "Done."
